$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Name:" cell to include the actual name
$ws.Range("A1").Value = "Name: Pandian"

# Update the "Git URL:" cell to include the actual repository URL
$ws.Range("A2").Value = "Git URL:https://github.com/Pandiyan927/react_gmail_assignment"

# Update the current selection to match the edited workbook (A2:B2, active cell A2)
$ws.Range("A2:B2").Select()
